# Fix for auto creation of FCF from KYCs
# Adds "Middle Name" and "Nationality" columns (Q & R) to the InvestorKyc sheet,
# and fills in sample values for the first two investor rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("Q1").Value = "Middle Name"
$ws.Range("R1").Value = "Nationality"

# New data for Investor 1
$ws.Range("Q2").Value = "Sir"
$ws.Range("R2").Value = "India"

# New data for Investor 2
$ws.Range("Q3").Value = "Sir"
$ws.Range("R3").Value = "US"

# Scroll / select the newly added column so it is visible, matching the
# author's final cursor position in the sheet.
$ws.Activate()
$excel.Goto($ws.Range("C1"), $true)
$ws.Range("R4").Select()
